# Weekly update: insert a new Haba (Macroferia Regional de Talca) price
# record as the new row 25, pushing the existing rows 25-39 down to 26-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 25 (shifts everything below it down by one,
# same as right-clicking row 25's header and choosing "Insert").
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new weekly price entry.
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44488
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 100112026
$ws.Range("G25").Value = "Haba"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 8000
$ws.Range("N25").Value = '$/saco 25 kilos'
$ws.Range("O25").Value = "Región de O'Higgins"
$ws.Range("P25").Value = 320
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
